# issue #5: stock data from json to db
#
# The stock worksheet (股票, sheet index 6) gains three new columns that
# mirror extra metadata the scraper now records alongside the existing
# property_category / date / legislator_name / legislator_id columns:
#   - "category"     inserted right before "date"
#   - "source_file"  appended after "legislator_id"
#   - "index"        appended after "source_file"
# A stray typo in one of the stock quantities ('vf9435 -> 9435) is also
# corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)   # 股票 (stock) sheet

$lastRow = 15                  # header is row 1, data rows are 2..15

# xlInsertFormatFromLeftOrAbove = 0 ; xlShiftToRight = -4161
$xlShiftToRight = -4161
$xlFormatFromLeftOrAbove = 0

# --- fix the stray typo in the quantity column for row 12 (聯格科技) ---
$ws.Cells.Item(12, 4).Value = "9435"

# --- insert a new "category" column before the existing "date" column ---
# Old layout:  H=property_category  I=date  J=legislator_name  K=legislator_id
# New layout:  H=property_category  I=category  J=date  K=legislator_name  L=legislator_id
$ws.Columns.Item(9).Insert($xlShiftToRight, $xlFormatFromLeftOrAbove)
$ws.Cells.Item(1, 9).Value = "category"

# --- append "source_file" right after legislator_id (column L / 12) ---
$ws.Columns.Item(13).Insert($xlShiftToRight, $xlFormatFromLeftOrAbove)
$ws.Cells.Item(1, 13).Value = "source_file"

# --- append "index" right after source_file (column M / 13) ---
$ws.Columns.Item(14).Insert($xlShiftToRight, $xlFormatFromLeftOrAbove)
$ws.Cells.Item(1, 14).Value = "index"

for ($row = 2; $row -le $lastRow; $row++) {
    # category column (I) -- matches the "normal" output folder
    $ws.Cells.Item($row, 9).Value = "normal"

    # source_file column (M) -- matches the tmp filename suffix
    $ws.Cells.Item($row, 13).Value = "tmpe4eb1"

    # index column (N) -- same as this row's original id in column A
    $idValue = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 14).Value2 = $idValue
}
